# Rename the existing "Sheet1" to "Mountains", add a new "Cities" sheet right
# after it, populate it with city data (closest cities referenced from the
# Mountains table, with population + altitude), and restyle it with a bold
# header row and thousands-separated numeric columns.

$wb = $excel.ActiveWorkbook

# --- Mountains (existing sheet) -------------------------------------------
$mountains = $wb.Worksheets.Item(1)
$mountains.Name = "Mountains"

# --- Cities (new sheet), inserted right after Mountains --------------------
$cities = $wb.Worksheets.Add($null, $mountains)
$cities.Name = "Cities"

# Header row
$cities.Range("A1").Value = "City"
$cities.Range("B1").Value = "Population"
$cities.Range("C1").Value = "Altitude"
$cities.Range("D1").Value = "Country"
$cities.Range("A1:D1").Font.Bold = $true

# Data rows: City, Population, Altitude (m), Country
$data = @(
  @("Kathmandu",  1442271, 1324, "Nepal"),
  @("Skardu",       26023, 2228, "Pakistan"),
  @("Darjeeling",  118805, 2045, "India"),
  @("Pokhara",     599504,  822, "Nepal"),
  @("Gorkha",      252201, 1130, "Nepal"),
  @("Gilgit",      216760, 1500, "Pakistan"),
  @("Joshimath",    48202, 1875, "India")
)

$r = 2
foreach ($row in $data) {
  $cities.Cells.Item($r, 1).Value = $row[0]
  $cities.Cells.Item($r, 2).Value = $row[1]
  $cities.Cells.Item($r, 3).Value = $row[2]
  $cities.Cells.Item($r, 4).Value = $row[3]
  $r++
}

# Thousands-separator formatting on Population/Altitude columns (C5 left
# unformatted, matching the source workbook).
$cities.Range("B2:B8").NumberFormat = "#,##0"
$cities.Range("C2:C4").NumberFormat = "#,##0"
$cities.Range("C6:C8").NumberFormat = "#,##0"

# --- Print setup on Mountains sheet ----------------------------------------
$mountains.PageSetup.PaperSize = 9
$mountains.PageSetup.Orientation = 1

# --- Selections / active sheet ---------------------------------------------
[void]$mountains.Range("D44").Select()
[void]$cities.Select()
[void]$cities.Range("G16").Select()
